$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's leftmost column (A) is blank; deleting it shifts the existing
# Year/CSO-002/CSO-004 table (previously occupying columns B:L) one column
# to the left (now A:K).
$ws.Range("A:A").Delete()

# Column A now holds the (previously unlabeled) year values. Give it a
# "Year" header, matching the bold/centered/wrapped formatting already used
# by the other header cells in row 1.
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("A1").Value = "Year"

# Reset the active selection to A1.
$ws.Range("A1").Select() | Out-Null
